$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "Para N = 100 e L = 0.5"
$ws.Range("A30").Value = "Classe"
$ws.Range("B30").Value = "Precisao"
$ws.Range("C30").Value = "Abrangencia"
$ws.Range("D30").Value = "F1"

$data = @(
    @(0, 0.994, 0.989, 0.992),
    @(1, 0.948, 0.995, 0.971),
    @(2, 0.978, 0.994, 0.986),
    @(3, 1, 0.951, 0.975),
    @(4, 0.973, 0.983, 0.978),
    @(5, 0.937, 0.984, 0.96),
    @(6, 1, 0.989, 0.994),
    @(7, 0.982, 0.922, 0.951),
    @(8, 0.975, 0.908, 0.94),
    @(9, 0.917, 0.978, 0.946)
)

$r = 31
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$ws.Range("A41").Value = "Media:"
$ws.Range("B41").Value = 0.97
$ws.Range("C41").Value = 0.969
$ws.Range("D41").Value = 0.969

$listObj = $ws.ListObjects.Add(1, $ws.Range("A30:D41"), 0, 1)
$listObj.Name = "Table1"
$listObj.TableStyle = "TableStyleMedium1"
